$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "UniqueId"
$ws.Range("F2").Value = "A23B58EC-8BB3-413F-BBBE-CCE71E470594"
$ws.Range("F3").Value = "6C5B3800-B097-4C52-85ED-A92321FCC9B0"

$ws.Range("F2").Select()
